# Insert two new rows before row 243, shifting the existing rows 243:370 down
# to 245:372 (matching the diff which adds two new data rows at the top of
# this date-block and pushes the rest of the table down by two rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(243).Resize(2).Insert()

# Fill in the new row 243 (first of the two inserted rows)
$ws.Range("A243").Value = 11
$ws.Range("B243").Value = "Vega Monumental Concepción"
$ws.Range("C243").Value = "Bíobío"
$ws.Range("D243").Value = 44813
$ws.Range("E243").Value = 8
$ws.Range("F243").Value = 100114014
$ws.Range("G243").Value = "Betarraga"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 600
$ws.Range("K243").Value = 700
$ws.Range("L243").Value = 800
$ws.Range("M243").Value = 750
$ws.Range("N243").Value = "$/paquete 5 unidades"
$ws.Range("O243").Value = "Región Metropolitana"
$ws.Range("P243").Value = 150
$ws.Range("Q243").Value = 5
$ws.Range("R243").Value = "Hortaliza"

# Fill in the new row 244 (second of the two inserted rows)
$ws.Range("A244").Value = 11
$ws.Range("B244").Value = "Vega Monumental Concepción"
$ws.Range("C244").Value = "Bíobío"
$ws.Range("D244").Value = 44813
$ws.Range("E244").Value = 8
$ws.Range("F244").Value = 100114014
$ws.Range("G244").Value = "Betarraga"
$ws.Range("H244").Value = "Sin especificar"
$ws.Range("I244").Value = "Segunda"
$ws.Range("J244").Value = 300
$ws.Range("K244").Value = 600
$ws.Range("L244").Value = 600
$ws.Range("M244").Value = 600
$ws.Range("N244").Value = "$/paquete 5 unidades"
$ws.Range("O244").Value = "Región Metropolitana"
$ws.Range("P244").Value = 120
$ws.Range("Q244").Value = 5
$ws.Range("R244").Value = "Hortaliza"
